$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.142.09"
$ws.Range("E2").Value = "  +6.16%  "
$ws.Range("D3").Value = "2.996.97"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'581.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.72%  "
$ws.Range("D6").Value = "'162.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.75%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D9").Value = "2.994.26"
$ws.Range("E9").Value = "  +3.70%  "
$ws.Range("D10").Value = "'6.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("E11").Value = "  +3.96%  "
$ws.Range("E12").Value = "  +5.91%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.14%  "
$ws.Range("D14").Value = "'34.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.80%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "66.159.14"
$ws.Range("E16").Value = "  +6.15%  "
$ws.Range("D17").Value = "3.495.11"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("D18").Value = "'6.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.14%  "
$ws.Range("D19").Value = "2.997.14"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").Value = "'452.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.21%  "
$ws.Range("D21").Value = "'13.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.50%  "
$ws.Range("D22").Value = "'0.684"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.29%  "
$ws.Range("D23").Value = "'7.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.27%  "
$ws.Range("D24").Value = "'82.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.92%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.18%  "
$ws.Range("D26").Value = "'12.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.89%  "
$ws.Range("D27").Value = "'10.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'8.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.88%  "
$ws.Range("E30").Value = "  +19.91%  "
$ws.Range("E31").Value = "  +5.89%  "
$ws.Range("D32").Value = "'0.0000103"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").Value = "'27.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.12%  "
$ws.Range("E34").Value = "  +5.25%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "'0.989"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.80%  "
$ws.Range("D37").Value = "'5.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.30%  "
$ws.Range("E38").Value = "  +9.07%  "
$ws.Range("D39").Value = "'49.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").Value = "'2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "'0.309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.57%  "
$ws.Range("D42").Value = "'44.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.08%  "
$ws.Range("E43").Value = "  +7.33%  "
$ws.Range("D44").Value = "'8.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.34%  "
$ws.Range("D45").Value = "'401.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.93%  "
$ws.Range("E46").Value = "  +6.53%  "
$ws.Range("D47").Value = "2.768.05"
$ws.Range("D48").Value = "'133.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D50").Value = "'23.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.73%  "
$ws.Range("E51").Value = "  +4.10%  "
